$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) New "Remarks" header in H2 (copy header format from B2:G2 style)
# ---------------------------------------------------------------
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H2").Value = "Remarks"

# ---------------------------------------------------------------
# 2) RHEL6+ block gets more results filled in.
#    Row 12 (RHEL6+, 1st data row): F12 -> Yes, G12 -> Yes (new, highlighted like end-of-block)
# ---------------------------------------------------------------
$ws.Range("F12").Value = "Yes"
$ws.Range("G5").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G12").Value = "Yes"

# Row 13 (2nd data row): F13 -> No, G13 -> No (no longer the highlighted last row),
# H13 gets a remark note.
$ws.Range("C4").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("H13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F13").Value = "No"
$ws.Range("G13").Value = "No"
$ws.Range("H13").Value = "fusermount becomes available only when the user is in the " + [char]0x2018 + "fuse" + [char]0x2019 + " group "

# Rows 14-16: three more fully-populated data rows (previously blank placeholder rows).
# Copy formatting from row 6, a plain data row with the same style pattern.
$ws.Range("A6:H6").Copy()
$ws.Range("A14:H14").PasteSpecial(-4122)
$ws.Range("A15:H15").PasteSpecial(-4122)
$ws.Range("A16:H16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B14").Value = "No"
$ws.Range("C14").Value = "No"
$ws.Range("D14").Value = "No"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "Yes"
$ws.Range("G14").Value = "No"

$ws.Range("B15").Value = "No"
$ws.Range("C15").Value = "No"
$ws.Range("D15").Value = "No"
$ws.Range("E15").Value = "No"
$ws.Range("F15").Value = "No"
$ws.Range("G15").Value = "No"

$ws.Range("B16").Value = "No"
$ws.Range("C16").Value = "No"
$ws.Range("D16").Value = "Yes"
$ws.Range("E16").Value = "No"
$ws.Range("F16").Value = "No"
$ws.Range("G16").Value = "No"

# ---------------------------------------------------------------
# 3) "Non-RHEL" section header moves from row 15 to row 18 (three rows down),
#    since it now comes after the 3 newly-filled rows above.
# ---------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A18").Value = "Non-RHEL"

# ---------------------------------------------------------------
# 4) Three new blank template rows appended at the bottom (32-34), matching
#    the existing blank rows just above them.
# ---------------------------------------------------------------
$ws.Range("A31:H31").Copy()
$ws.Range("A32:H32").PasteSpecial(-4122)
$ws.Range("A33:H33").PasteSpecial(-4122)
$ws.Range("A34:H34").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows("32:34").RowHeight = $ws.Rows("31").RowHeight

# ---------------------------------------------------------------
# 5) Column H is widened to fit the new Remarks text.
# ---------------------------------------------------------------
$ws.Columns("H").ColumnWidth = 55.7969

# ---------------------------------------------------------------
# 6) The Yes/No data-validation dropdown now covers the extended data
#    range B3:G18 (was B3:G15).
# ---------------------------------------------------------------
$ws.Range("B3:G18").Validation.Delete()
$ws.Range("B3:G18").Validation.Add(3, 1, 1, ",Yes,No")
